# LOM3203.xlsx update
# The row that only held "5982760 - Carlos Alberto Baldan" (row 13, no label in
# column A) is removed entirely; every row below it shifts up by one. On top
# of that shift, several B/C cells get new text (the professor name and a
# handful of long paragraphs are swapped around / replaced).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the old row 13 ("5982760 - Carlos Alberto Baldan" only, no A label).
#    Excel shifts rows 14-25 up to become rows 13-24, carrying their row
#    heights and existing content along for free.
$ws.Rows("13").Delete()

# 2) Patch the handful of cells whose text content changed, using the
#    post-shift row numbers.

# Row 10 ("Objetivos:") - long objectives paragraph -> professor name
$ws.Range("B10").Value = "5982760 - Carlos Alberto Baldan"
$ws.Range("C10").Value = "5982760 - Carlos Alberto Baldan"

# Row 13 ("Programa resumido:", was row 14) -> "Semestral"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 ("Programa:", was row 16) -> "01/01/2012"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"

# Row 18 ("Método:", was row 19) -> professor name
$ws.Range("B18").Value = "5982760 - Carlos Alberto Baldan"
$ws.Range("C18").Value = "5982760 - Carlos Alberto Baldan"

# Row 19 ("Critério:", was row 20) -> the "curso desenvolvido..." text
$ws.Range("B19").Value = "O curso é desenvolvido através de aulas expositivas e práticas em laboratório."
$ws.Range("C19").Value = "O curso é desenvolvido através de aulas expositivas e práticas em laboratório."

# Row 20 ("Norma de recuperação:", was row 21) -> the "média aritmética..." text
$ws.Range("B20").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."
$ws.Range("C20").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."

# Row 21 ("Bibliografia:", was row 22) -> the "aplicação de uma prova escrita..." text
$ws.Range("B21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
